# Split the single run containing
#   "QUEDARÁ FACULTADO PARA DISPONER LIBREMENTE DEL BIEN INMUEBLE..."
# into three runs:
#   "QUEDARÁ FACULTAD" / "{{SEXO_14}}" / " PARA DISPONER LIBREMENTE..."
# (i.e. turn the literal "FACULTADO" into "FACULTAD" + a {{SEXO_14}} merge field)

$d = $word.ActiveDocument

# Locate the target sentence and remember its bounds.
$target = $d.Content
$found = $target.Find.Execute("QUEDARÁ FACULTADO")
if (-not $found) {
    throw "Could not find 'QUEDARÁ FACULTADO' in the document"
}
$sentenceStart = $target.Start
$sentenceEnd = $target.End

# The trailing "O" of "FACULTADO" is the last character of the match.
$oStart = $sentenceEnd - 1
$oEnd = $sentenceEnd

# Split that trailing "O" into its own run first (formatting toggle trick,
# see below) so the later text replacement stays scoped to a single run
# instead of being absorbed into a neighboring run.
$rO = $d.Range($oStart, $oEnd)
$rO.Font.Bold = $true
$rO.Font.Bold = $false

# Replace that "O" with the merge-field placeholder text, scoped to just
# that run via Find/Replace (keeps the run's original, bare <w:r> -- a
# plain Range.Text assignment instead tends to get absorbed into the
# preceding run and drags its rsid attributes along with it).
$rO.Find.Execute("O", $true, $false, $false, $false, $false, $true, 1, `
                  $false, "{{SEXO_14}}", 2) | Out-Null

# Re-find the stable anchor text now that the placeholder is in place.
$verify = $d.Content
$verify.Find.Execute("QUEDARÁ FACULTAD{{SEXO_14}}")
$newStart = $verify.Start

$placeholderStart = $newStart + 16
$placeholderEnd = $placeholderStart + 11  # length of "{{SEXO_14}}"

# Now force the three runs apart by toggling a character-formatting
# property on and back off at each boundary; the engine keeps runs split
# at a formatting boundary even after the property reverts to its
# original value, which is exactly how Word itself leaves behind
# independent <w:r> elements with identical <w:rPr>.
$rFirst = $d.Range($newStart, $placeholderStart)
$rFirst.Font.Bold = $true
$rFirst.Font.Bold = $false

$rPlaceholder = $d.Range($placeholderStart, $placeholderEnd)
$rPlaceholder.Font.Bold = $true
$rPlaceholder.Font.Bold = $false

Write-Output "QUEDARÁ FACULTAD / {{SEXO_14}} split complete"
